# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.185.14'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.908.87'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8301'
$ws.Range('E5').Value = '  +4.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.05'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3277'
$ws.Range('E8').Value = '  +3.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.91'
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07070'
$ws.Range('E10').Value = '  +1.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08104'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7671'
$ws.Range('E12').Value = '  +2.15%  '
$ws.Range('D13').Value = '1.917.31'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.90'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').Value = '30.190.65'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.22'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.911'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.90'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007794'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = '2.160.89'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.051'
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1705'
$ws.Range('E25').Value = '  +22.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.332'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.12'
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.106'
$ws.Range('E29').Value = '  +2.32%  '
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.06102'
$ws.Range('E31').Value = '  +7.45%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.530'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.309'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.091'
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.276'
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7341'
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.713'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01938'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.796'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4471'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.33'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.957'
$ws.Range('E42').Value = '  -3.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8578'
$ws.Range('E43').Value = '  +2.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.910'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.45'
$ws.Range('E46').Value = '  +1.34%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.010.28'
$ws.Range('E47').Value = '  +2.37%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.592'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.869'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').Value = '2.063.63'
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.551'
$ws.Range('E51').Value = '  +3.20%  '
